# "Added gifs to BPL" -- trims the trailing tab run in the "Valuation - DCF"
# title placeholder (on both the slide that already carried the title and
# its duplicate) from 7 tabs down to 4, leaving the rest of each run's
# formatting/text untouched.

$p = $ppt.ActivePresentation

$newLead = "Valuation - DCF " + "`t`t`t`t"

foreach ($idx in 28, 29) {
    $slide = $p.Slides.Item($idx)
    $title = $slide.Shapes.Item(1)
    $full = $title.TextFrame.TextRange
    # First run is "Valuation - DCF " followed by 7 tabs (23 chars total);
    # replace exactly that run's characters so its own run/formatting is
    # preserved and only the <a:t> payload shrinks to 4 trailing tabs.
    $lead = $full.Characters(1, 23)
    $lead.Text = $newLead
}
